$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17/18: Litecoin / ShibaInu swap ---
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'87.10"
$ws.Range("E17").Value = "  +1.46%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.00001045"
$ws.Range("E18").Value = "  +1.55%  "

# --- Price / Volume updates ---
$ws.Range("D2").Value = "27.783.33"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "1.877.43"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'331.62"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4723"
$ws.Range("E7").Value = "  +4.33%  "
$ws.Range("D8").Value = "'0.3950"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "'47.96"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "'0.08078"
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").Value = "'1.029"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("E12").Value = "  +3.79%  "
$ws.Range("D13").Value = "1.870.48"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'5.968"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "'7.137"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D19").Value = "'0.06659"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("D20").Value = "'17.19"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "27.785.60"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").Value = "'5.528"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").Value = "'2.303"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "2.099.82"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "'159.06"
$ws.Range("E27").Value = "  +3.38%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'2.109"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").Value = "'5.595"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("D31").Value = "'122.36"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").Value = "'0.9862"
$ws.Range("E32").Value = "  +5.48%  "
$ws.Range("D33").Value = "'0.09549"
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").Value = "'3.591"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'5.344"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D37").Value = "'0.06116"
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("D38").Value = "'0.02258"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").Value = "'1.232"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "'8.165"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").Value = "'0.6039"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").Value = "'0.1907"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'10.26"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").Value = "'12.17"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").Value = "'1.949"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").Value = "'3.378"
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "'0.06896"
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("D50").Value = "'114.83"
$ws.Range("E50").Value = "  +5.84%  "
$ws.Range("E51").Value = "  +1.96%  "
